$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the publication detail columns (authors, year, title, journal, abstract, date, link)
# for the two data rows, keeping only the short_id column (H) intact.
$ws.Range("A2:G2").ClearContents()
$ws.Range("A3:G3").ClearContents()

# Reset the row heights back to the default (they were manually stretched to
# fit the long abstract text previously stored in column E).
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# Update the selection / active cell to match the new focus point.
$ws.Range("B3").Select()
